$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 17.39238349376416
$ws.Range("C2").Value = 5.26521440389474
$ws.Range("D2").Value = 11.08490074222781
$ws.Range("E2").Value = 11.15595206546577
$ws.Range("F2").Value = 55.9786459778256
$ws.Range("K2").Value = 13.68452227586697
$ws.Range("L2").Value = 10.21633361865078
$ws.Range("B3").Value = 17.34390197730502
$ws.Range("C3").Value = 5.180421910722615
$ws.Range("D3").Value = 10.94327030893418
$ws.Range("E3").Value = 11.13865921087243
$ws.Range("F3").Value = 54.83860886638152
$ws.Range("K3").Value = 13.68601242488004
$ws.Range("L3").Value = 10.21941548474439
$ws.Range("B4").Value = 17.32078501600713
$ws.Range("C4").Value = 5.126050064254225
$ws.Range("D4").Value = 10.85431888403606
$ws.Range("E4").Value = 11.12945863325354
$ws.Range("F4").Value = 54.12850773851144
$ws.Range("K4").Value = 13.69238407058905
$ws.Range("L4").Value = 10.22337281166613
$ws.Range("B5").Value = 17.31304411219203
$ws.Range("C5").Value = 5.103313453015731
$ws.Range("D5").Value = 10.81758476478675
$ws.Range("E5").Value = 11.12606712026626
$ws.Range("F5").Value = 53.83687583918103
$ws.Range("K5").Value = 13.69635027585462
$ws.Range("L5").Value = 10.22550463785311
$ws.Range("B6").Value = 17.31186033411974
$ws.Range("C6").Value = 5.099503077952554
$ws.Range("D6").Value = 10.81145619092508
$ws.Range("E6").Value = 11.12552560927629
$ws.Range("F6").Value = 53.78832255196463
$ws.Range("K6").Value = 13.69709149257987
$ws.Range("L6").Value = 10.22588997657659
$ws.Range("B7").Value = 17.32067381202957
$ws.Range("C7").Value = 5.125745775879287
$ws.Range("D7").Value = 10.85382542177844
$ws.Range("E7").Value = 11.12941144369368
$ws.Range("F7").Value = 54.12458345809926
$ws.Range("K7").Value = 13.69243201855069
$ws.Range("L7").Value = 10.22339946035141
$ws.Range("B8").Value = 17.37429161009148
$ws.Range("C8").Value = 5.236457678354028
$ws.Range("D8").Value = 11.03648823692742
$ws.Range("E8").Value = 11.14969622488079
$ws.Range("F8").Value = 55.58782715190901
$ws.Range("K8").Value = 13.68390247110765
$ws.Range("L8").Value = 10.21696765057388
$ws.Range("B9").Value = 17.53178714232538
$ws.Range("C9").Value = 5.435214008443507
$ws.Range("D9").Value = 11.37827846676463
$ws.Range("E9").Value = 11.20065213200483
$ws.Range("F9").Value = 58.3650935013718
$ws.Range("K9").Value = 13.71052145337238
$ws.Range("L9").Value = 10.22073673742895
$ws.Range("B10").Value = 17.67869170758803
$ws.Range("C10").Value = 5.57003335069198
$ws.Range("D10").Value = 11.61858511635892
$ws.Range("E10").Value = 11.24480284124088
$ws.Range("F10").Value = 60.33429461897306
$ws.Range("K10").Value = 13.7564834454268
$ws.Range("L10").Value = 10.233477230326
$ws.Range("B11").Value = 17.75209437277612
$ws.Range("C11").Value = 5.628931918819641
$ws.Range("D11").Value = 11.72539848105151
$ws.Range("E11").Value = 11.26632042125926
$ws.Range("F11").Value = 61.21173445020488
$ws.Range("K11").Value = 13.78309369605076
$ws.Range("L11").Value = 10.24143052699367
$ws.Range("B12").Value = 17.78081536612217
$ws.Range("C12").Value = 5.650885577634046
$ws.Range("D12").Value = 11.76547277334662
$ws.Range("E12").Value = 11.27467211172668
$ws.Range("F12").Value = 61.54113215065977
$ws.Range("K12").Value = 13.79398526552772
$ws.Range("L12").Value = 10.24475138932839
$ws.Range("B13").Value = 17.77458899381449
$ws.Range("C13").Value = 5.646173031961275
$ws.Range("D13").Value = 11.75685888559355
$ws.Range("E13").Value = 11.2728644214339
$ws.Range("F13").Value = 61.47032159145424
$ws.Range("K13").Value = 13.79160342619226
$ws.Range("L13").Value = 10.24402245456677
$ws.Range("B14").Value = 17.75443885218387
$ws.Range("C14").Value = 5.63074507433194
$ws.Range("D14").Value = 11.72870298118672
$ws.Range("E14").Value = 11.26700345825355
$ws.Range("F14").Value = 61.2388929730492
$ws.Range("K14").Value = 13.7839734553056
$ws.Range("L14").Value = 10.24169754843486
$ws.Range("B15").Value = 17.74221612893093
$ws.Range("C15").Value = 5.621249424181739
$ws.Range("D15").Value = 11.71140762566944
$ws.Range("E15").Value = 11.26343986247935
$ws.Range("F15").Value = 61.09675597878292
$ws.Range("K15").Value = 13.77940582044155
$ws.Range("L15").Value = 10.24031369526088
$ws.Range("B16").Value = 17.67402527167921
$ws.Range("C16").Value = 5.566135233863638
$ws.Range("D16").Value = 11.61155306125455
$ws.Range("E16").Value = 11.24342523505672
$ws.Range("F16").Value = 60.27656245110121
$ws.Range("K16").Value = 13.75485873689943
$ws.Range("L16").Value = 10.23300078806658
$ws.Range("B17").Value = 17.63386247659811
$ws.Range("C17").Value = 5.53170144436405
$ws.Range("D17").Value = 11.5496444962412
$ws.Range("E17").Value = 11.23151216919274
$ws.Range("F17").Value = 59.76852985596768
$ws.Range("K17").Value = 13.74125711328767
$ws.Range("L17").Value = 10.22906643214138
$ws.Range("B18").Value = 17.61138213081494
$ws.Range("C18").Value = 5.511667190922637
$ws.Range("D18").Value = 11.51380148566066
$ws.Range("E18").Value = 11.22479515931541
$ws.Range("F18").Value = 59.47460764524099
$ws.Range("K18").Value = 13.73397090240634
$ws.Range("L18").Value = 10.22700665973052
$ws.Range("B19").Value = 17.60387778875282
$ws.Range("C19").Value = 5.504844644128223
$ws.Range("D19").Value = 11.50162574763271
$ws.Range("E19").Value = 11.2225441709555
$ws.Range("F19").Value = 59.37480327075571
$ws.Range("K19").Value = 13.73159629339984
$ws.Range("L19").Value = 10.2263441806814
$ws.Range("B20").Value = 17.63807382518252
$ws.Range("C20").Value = 5.535390667738597
$ws.Range("D20").Value = 11.55625917422675
$ws.Range("E20").Value = 11.23276637462406
$ws.Range("F20").Value = 59.82279003051624
$ws.Range("K20").Value = 13.74264947572719
$ws.Range("L20").Value = 10.22946423171401
$ws.Range("B21").Value = 17.76033250605822
$ws.Range("C21").Value = 5.635286138782448
$ws.Range("D21").Value = 11.73698329863954
$ws.Range("E21").Value = 11.26871946567377
$ws.Range("F21").Value = 61.30694880769468
$ws.Range("K21").Value = 13.78619249534293
$ws.Range("L21").Value = 10.24237205048774
$ws.Range("B22").Value = 17.8456163570652
$ws.Range("C22").Value = 5.698533784816609
$ws.Range("D22").Value = 11.85291331348867
$ws.Range("E22").Value = 11.29340138786529
$ws.Range("F22").Value = 62.26010946330663
$ws.Range("K22").Value = 13.8193965974034
$ws.Range("L22").Value = 10.2526091581951
$ws.Range("B23").Value = 17.79961369054638
$ws.Range("C23").Value = 5.664964034080498
$ws.Range("D23").Value = 11.79124345780891
$ws.Range("E23").Value = 11.28012072316012
$ws.Range("F23").Value = 61.75300098460267
$ws.Range("K23").Value = 13.8012426954007
$ws.Range("L23").Value = 10.24698105692527
$ws.Range("B24").Value = 17.63616797362892
$ws.Range("C24").Value = 5.533723510547295
$ws.Range("D24").Value = 11.55326945892879
$ws.Range("E24").Value = 11.23219893709209
$ws.Range("F24").Value = 59.79826474258444
$ws.Range("K24").Value = 13.74201832643
$ws.Range("L24").Value = 10.22928375681911
$ws.Range("B25").Value = 17.4836413287916
$ws.Range("C25").Value = 5.383412313841299
$ws.Range("D25").Value = 11.2876563393874
$ws.Range("E25").Value = 11.18568015730479
$ws.Range("F25").Value = 57.62524153061084
$ws.Range("K25").Value = 13.6986765167934
$ws.Range("L25").Value = 10.21796399885296
